# Auto-generated edit script applying cryptos list update
# (GitHub Actions scrape refresh, Wed Oct 30 02:57:22 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '72.365.93'
Set-TextCell 'E2' '  +1.56%  '
Set-TextCell 'D3' '2.628.88'
Set-TextCell 'E3' '  +0.58%  '
Set-TextCell 'D4' '1.00'
Set-TextCell 'E4' '  +0.04%  '
Set-TextCell 'D5' '601.21'
Set-TextCell 'E5' '  -0.99%  '
Set-TextCell 'D6' '179.63'
Set-TextCell 'E6' '  -0.76%  '
Set-TextCell 'E7' '  -0.02%  '
Set-TextCell 'D8' '0.525'
Set-TextCell 'E8' '  +0.05%  '
Set-TextCell 'D9' '0.175'
Set-TextCell 'E9' '  +5.01%  '
Set-TextCell 'D10' '2.628.07'
Set-TextCell 'E10' '  +0.63%  '
Set-TextCell 'E11' '  +1.23%  '
Set-TextCell 'D12' '0.360'
Set-TextCell 'E12' '  +3.06%  '
Set-TextCell 'D13' '5.03'
Set-TextCell 'E13' '  -0.61%  '
Set-TextCell 'B14' 'WrappedliquidstakedEther2.0'
Set-TextCell 'C14' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 'D14' '3.134.13'
Set-TextCell 'E14' '  +0.63%  '
Set-TextCell 'B15' 'ShibaInu'
Set-TextCell 'C15' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 'D15' '0.0000190'
Set-TextCell 'E15' '  +2.34%  '
Set-TextCell 'D16' '72.300.17'
Set-TextCell 'E16' '  +1.54%  '
Set-TextCell 'D17' '26.62'
Set-TextCell 'E17' '  -0.82%  '
Set-TextCell 'D18' '2.629.43'
Set-TextCell 'E18' '  +1.61%  '
Set-TextCell 'D19' '11.87'
Set-TextCell 'E19' '  +4.14%  '
Set-TextCell 'D20' '379.95'
Set-TextCell 'E20' '  +2.47%  '
Set-TextCell 'D21' '7.87'
Set-TextCell 'E21' '  +0.13%  '
Set-TextCell 'E22' '  -0.28%  '
Set-TextCell 'D23' '2.08'
Set-TextCell 'E23' '  +9.58%  '
Set-TextCell 'D24' '73.30'
Set-TextCell 'E24' '  +1.54%  '
Set-TextCell 'E25' '  +0.07%  '
Set-TextCell 'E26' '  -1.40%  '
Set-TextCell 'D27' '10.11'
Set-TextCell 'E27' '  +4.78%  '
Set-TextCell 'D28' '2.766.35'
Set-TextCell 'E28' '  +0.71%  '
Set-TextCell 'D29' '1.00'
Set-TextCell 'E29' '  +0.08%  '
Set-TextCell 'D30' '0.0₃0953'
Set-TextCell 'E30' '  -0.26%  '
Set-TextCell 'E31' '  +0.27%  '
Set-TextCell 'D32' '519.38'
Set-TextCell 'E32' '  -2.79%  '
Set-TextCell 'D33' '1.31'
Set-TextCell 'E33' '  -0.77%  '
Set-TextCell 'E34' '  -0.85%  '
Set-TextCell 'D35' '0.999'
Set-TextCell 'E35' '  -0.11%  '
Set-TextCell 'D36' '164.58'
Set-TextCell 'E36' '  +0.09%  '
Set-TextCell 'D37' '19.31'
Set-TextCell 'E37' '  +0.32%  '
Set-TextCell 'E38' '  -6.31%  '
Set-TextCell 'D39' '19.09'
Set-TextCell 'E39' '  +0.65%  '
Set-TextCell 'E40' '  +1.12%  '
Set-TextCell 'D41' '1.85'
Set-TextCell 'E41' '  +1.56%  '
Set-TextCell 'B42' 'dogwifhat'
Set-TextCell 'C42' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D42' '2.64'
Set-TextCell 'E42' '  +2.27%  '
Set-TextCell 'B43' 'RenderToken'
Set-TextCell 'C43' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell 'D43' '5.07'
Set-TextCell 'E43' '  -0.05%  '
Set-TextCell 'B44' 'USDe'
Set-TextCell 'C44' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell 'D44' '1.00'
Set-TextCell 'E44' '  -0.15%  '
Set-TextCell 'E45' '  +0.89%  '
Set-TextCell 'D46' '39.32'
Set-TextCell 'E46' '  -1.13%  '
Set-TextCell 'D47' '150.84'
Set-TextCell 'E47' '  -2.48%  '
Set-TextCell 'E48' '  +0.57%  '
Set-TextCell 'E49' '  +2.09%  '
Set-TextCell 'E50' '  +1.25%  '
Set-TextCell 'D51' '0.0₆0261'
Set-TextCell 'E51' '  -3.37%  '
